$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated financial figures
$ws.Range("D2").Value = 1664
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = -12
$ws.Range("H2").Value = -16
$ws.Range("I2").Value = -16
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2554
$ws.Range("L2").Value = 1290
$ws.Range("M2").Value = 1263
$ws.Range("N2").Value = 1263
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 164
$ws.Range("Q2").Value = -31
$ws.Range("R2").Value = -54
$ws.Range("S2").Value = 27
$ws.Range("T2").Value = 55
$ws.Range("U2").Value = -85
$ws.Range("V2").Value = 928
$ws.Range("W2").Value = 1.78
$ws.Range("X2").Value = -0.96
$ws.Range("Y2").Value = -1.22
$ws.Range("Z2").Value = -0.62
$ws.Range("AA2").Value = 102.11
$ws.Range("AB2").Value = 663.23
$ws.Range("AC2").Value = -104
$ws.Range("AD2").Value = -41.61
$ws.Range("AE2").Value = 8292
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 3.48
$ws.Range("AI2").Value = -144.83
$ws.Range("AJ2").Value = 15225000

# Row 3: updated financial figures
$ws.Range("D3").Value = 1677
$ws.Range("E3").Value = -23
$ws.Range("F3").Value = -23
$ws.Range("G3").Value = -30
$ws.Range("H3").Value = -37
$ws.Range("I3").Value = -37
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2335
$ws.Range("L3").Value = 1150
$ws.Range("M3").Value = 1186
$ws.Range("N3").Value = 1185
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 164
$ws.Range("Q3").Value = 42
$ws.Range("R3").Value = 151
$ws.Range("S3").Value = -165
$ws.Range("T3").Value = 45
$ws.Range("U3").Value = -3
$ws.Range("V3").Value = 790
$ws.Range("W3").Value = -1.37
$ws.Range("X3").Value = -2.22
$ws.Range("Y3").Value = -3.04
$ws.Range("Z3").Value = -1.52
$ws.Range("AA3").Value = 96.96
$ws.Range("AB3").Value = 617.06
$ws.Range("AC3").Value = -245
$ws.Range("AD3").Value = -20.76
$ws.Range("AE3").Value = 7781
$ws.Range("AF3").Value = 0.65
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 2.95
$ws.Range("AI3").Value = -59.8
$ws.Range("AJ3").Value = 15225000

# Row 4: updated financial figures
$ws.Range("D4").Value = 1744
$ws.Range("E4").Value = 43
$ws.Range("F4").Value = 43
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2307
$ws.Range("L4").Value = 1141
$ws.Range("M4").Value = 1166
$ws.Range("N4").Value = 1165
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 164
$ws.Range("Q4").Value = 133
$ws.Range("R4").Value = -59
$ws.Range("S4").Value = -42
$ws.Range("T4").Value = 49
$ws.Range("U4").Value = 84
$ws.Range("V4").Value = 767
$ws.Range("W4").Value = 2.48
$ws.Range("X4").Value = 0.19
$ws.Range("Y4").Value = 0.28
$ws.Range("Z4").Value = 0.15
$ws.Range("AA4").Value = 97.87
$ws.Range("AB4").Value = 606.08
$ws.Range("AC4").Value = 21
$ws.Range("AD4").Value = 285.21
$ws.Range("AE4").Value = 7649
$ws.Range("AF4").Value = 0.8
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 2.46
$ws.Range("AI4").Value = 701.33
$ws.Range("AJ4").Value = 15225000

# Row 5: updated financial figures
$ws.Range("D5").Value = 1972
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = -48
$ws.Range("H5").Value = -43
$ws.Range("I5").Value = -43
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2293
$ws.Range("L5").Value = 1185
$ws.Range("M5").Value = 1107
$ws.Range("N5").Value = 1106
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 164
$ws.Range("Q5").Value = 41
$ws.Range("R5").Value = 140
$ws.Range("S5").Value = -138
$ws.Range("T5").Value = 40
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 640
$ws.Range("W5").Value = 1.01
$ws.Range("X5").Value = -2.18
$ws.Range("Y5").Value = -3.82
$ws.Range("Z5").Value = -1.87
$ws.Range("AA5").Value = 107.03
$ws.Range("AB5").Value = 563.57
$ws.Range("AC5").Value = -285
$ws.Range("AD5").Value = -19.3
$ws.Range("AE5").Value = 7264
$ws.Range("AF5").Value = 0.76
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 2.73
$ws.Range("AI5").Value = -52.64
$ws.Range("AJ5").Value = 15225000

# Row 6: updated financial figures
$ws.Range("D6").Value = 2013
$ws.Range("E6").Value = -10
$ws.Range("F6").Value = -10
$ws.Range("G6").Value = 59
$ws.Range("H6").Value = -12
$ws.Range("I6").Value = -11
$ws.Range("K6").Value = 1993
$ws.Range("L6").Value = 925
$ws.Range("M6").Value = 1067
$ws.Range("N6").Value = 1064
$ws.Range("P6").Value = 164
$ws.Range("Q6").Value = -116
$ws.Range("R6").Value = 306
$ws.Range("S6").Value = -166
$ws.Range("T6").Value = 33
$ws.Range("U6").Value = -149
$ws.Range("V6").Value = 508
$ws.Range("W6").Value = -0.52
$ws.Range("X6").Value = -0.61
$ws.Range("Y6").Value = -1.04
$ws.Range("Z6").Value = -0.57
$ws.Range("AA6").Value = 86.71
$ws.Range("AB6").Value = 540.54
$ws.Range("AC6").Value = -74
$ws.Range("AD6").Value = -56.41
$ws.Range("AE6").Value = 6988
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 3.59
$ws.Range("AI6").Value = -202.43
$ws.Range("AJ6").Value = 15225000

# Rows 7-9: data no longer available for these periods, remove values
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
